# Update the "F" column (想去人数 / want-to-go count) values on the
# "展览" and "全部类型" worksheets, per the commit's regenerated data output.
# Row numbers diverge between the two sheets after row 29 because
# "全部类型" interleaves rows from the other category sheets, so each
# sheet gets its own row->value map.

$wb = $excel.ActiveWorkbook

$updatesExhibition = [ordered]@{
    2  = 186
    4  = 142
    5  = 1297
    6  = 18141
    7  = 363
    8  = 260
    10 = 6843
    11 = 687
    12 = 159
    13 = 14
    15 = 63
    17 = 155
    18 = 1300
    19 = 223
    21 = 655
    25 = 273
    26 = 985
    27 = 125
    29 = 534
    30 = 33
    32 = 71
    33 = 12068
    34 = 1281
    37 = 280
    38 = 3919
    39 = 301
}

$updatesAllTypes = [ordered]@{
    2  = 186
    4  = 142
    5  = 1297
    6  = 18141
    7  = 363
    8  = 260
    10 = 6843
    11 = 687
    12 = 159
    13 = 14
    15 = 63
    17 = 155
    18 = 1300
    19 = 223
    21 = 655
    25 = 273
    26 = 985
    27 = 125
    29 = 534
    32 = 33
    34 = 71
    35 = 12068
    36 = 1281
    39 = 280
    40 = 3919
    41 = 301
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $ws1.Range("F$row").Value = $updatesExhibition[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $ws4.Range("F$row").Value = $updatesAllTypes[$row]
}
